$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2025-07-27 00:00:00"
$ws.Range("G2").Value = "• Singapore's financial job market is slowing due to higher labor costs, tighter visa requirements, and the rise of AI, leading companies to expand in neighboring countries like Malaysia and offshore roles to India.`n• DBS plans to cut 4,000 temporary jobs by using AI to replace `"project`" roles, and the Singaporean government is collaborating with six banks to identify jobs impacted by AI and upskill employees.`n• Hong Kong's financial job market is experiencing a resurgence driven by capital inflows from China, record-high IPOs, and its focus on front-office activities like sales and trading, which are less susceptible to automation compared to Singapore's back and middle office roles."
$ws.Range("G3").Value = "• Manulife and RBC are using internal tests to evaluate AI models, finding that standard third-party benchmarks don't accurately reflect performance for their specific business needs. Manulife's global chief AI officer, Jodie Wallis, stated, `"The benchmarks can actually lead you astray.`"`n• Manulife has developed its own testing program based on 47 real-world use cases, cutting evaluation time from weeks to minutes, while RBC is standardizing its AI platform using Cohere's technology, prioritizing security and scalability over raw AI power.`n• RBC is co-developing a secure version of Cohere's North product to meet its regulatory requirements and enable AI application development with sensitive data, aiming to generate up to `$1 billion in earnings using AI by 2027."
$ws.Range("G4").Value = "• Global Finance has introduced its first AI in Finance awards, acknowledging the increasing impact of AI, especially generative AI (GenAI), on the finance sector.`n• Banks are taking a risk-based, governance-first approach to GenAI, focusing on legal, regulatory, and ethical considerations to ensure responsible, secure, and ethical deployment, as exemplified by DBS Bank's PURE Framework for Data Ethics.`n• Banks are struggling with calculating the ROI of AI projects, which poses a challenge for securing funding for experimental and innovative AI initiatives beyond simple automation, according to Andy Schmidt, vice president and global industry lead for banking at CGI."
$ws.Range("G5").Value = "• Banks like JPMorgan Chase, Capital One, and BNY Mellon are starting to adopt agentic AI, which utilizes large language models to perform tasks with minimal human intervention. However, Gartner predicts that 40% of agentic AI projects will be shut down by 2027 due to escalating costs, unclear business value, or inadequate risk controls.`n• Agentic AI use cases in banking include performing `"four eyes checks`" and `"six eyes checks`" for high-value transactions and verifying the accuracy of data captured by traders on trading floors, mimicking middle office functions. Some banks are framing agentic AI as assistants to employees, while others are considering them as autonomous `"digital workers`" with the ability to sense, perceive, and adapt.`n• Challenges in managing AI agents include handling access entitlements, redesigning processes to suit AI capabilities, and continuously monitoring AI performance for model drift. Best practices include providing frequent feedback, involving AI creators in management, and training employees to effectively collaborate with AI agents."
$ws.Range("G6").Value = "• Standard Chartered is comfortable being a `"good follower`" in AI adoption, prioritizing risk mitigation over being at the `"cutting edge of innovation,`" given the potential for significant repercussions from AI-related errors.`n• Ant International highlights key AI risks: hallucinations, AI agent vulnerabilities to external attacks, and deepfakes, including the potential for fake customers used as attack vectors, while also noting AI's ability to enhance entry-level financial investigators' capabilities.`n• BNP Paribas Wealth Management observes a generational divide in AI trust: younger clients embrace AI for its speed and transparency in thematic investing, while older clients view it cautiously, as a `"support tool, and not something to be used for investing.`"`n• Larger, established banks are now driving the AI agenda, a reversal from previous digital transformations where smaller companies led the way, while startups face challenges like long-term horizons and complex documentation required by larger institutions.`n• Amber Group is exploring `"AgentFi,`" finance driven by autonomous AI agents, believing crypto can give AI agents `"financial freedom`" and the resources to execute their decisions, demonstrated by their `"AgentFi Ambassador`" AI, Mia, currently managing liquidity."
$ws.Range("G7").Value = "• Citi's M&A revenues increased 52% year-on-year in Q2, while JPMorgan's rose only 8%, with Citi ranking fourth globally for M&A in the first half of 2025.`n• Citi has implemented Gen AI tools and automated 740,000 code reviews, saving 100,000 hours of developer time, but this has led to severance charges related to the `"realignment of the technology workforce.`"`n• Citi's DCM revenues fell by 12% in Q2, even while they rose by a similar proportion at JPMorgan."
$ws.Range("G8").Value = "• Aditya Bhasin, Bank of America's chief technology and information officer, is retiring and will advise/invest in mid-sized companies.`n• Hari Gopalkrishnan, formerly in senior tech roles at Citigroup, will replace Bhasin, having joined BofA in 2011.`n• Bank of America has invested heavily in technology, resulting in over 59 million digital users, 7,800+ patent filings, and 3 billion client interactions via its AI virtual assistant, Erica, since 2018."
$ws.Range("G9").Value = "• Goldman Sachs is testing Devin, an autonomous software engineer from AI startup Cognition, with plans to integrate it into its team of 12,000 developers.`n• Goldman Sachs tech chief Marco Argenti envisions scaling to `"hundreds`" or potentially `"thousands`" of Devins, anticipating a 3-4x increase in worker productivity compared to previous AI tools.`n• Devin will focus on tasks like updating internal code to newer programming languages, effectively handling `"drudgery`" work under human supervision."
$ws.Range("G10").Value = "• HSBC UK redesigned its mobile banking app, its first major overhaul since 2012, to address outdated design, complex navigation, and technical constraints that hindered speed of change and adaptability. The redesign took 18 months.`n• HSBC is leveraging generative AI and Google Dialogflow to enhance its conversational banking capabilities and digital journeys, `"from design and improving engineering productivity and accuracy,`" alongside its natural language understanding.`n• The new app, initially launched in the UK, features improved discoverability of products, a better chat experience, and a loosely coupled, event-driven architecture for greater flexibility; future enhancements include in-app credit card reward redemption and a fraud management hub."
$ws.Range("G11").Value = "• UnlikelyAI, founded by the creator of Evi (acquired by Amazon and integrated into Alexa), uses a neurosymbolic approach to combine neural networks and symbolic systems, aiming to address unpredictability issues in AI for financial services.`n• Lloyds Banking Group is testing UnlikelyAI's platform within its innovation sandbox to evaluate if the technology can maintain transparency and explainability standards while delivering operational benefits, aligning with banking sector requirements.`n• UnlikelyAI's technology focuses on reliability and compliance, enabling banks to deploy conversational AI while meeting financial services operational standards, potentially overcoming challenges that have limited AI adoption in customer-facing banking applications."
$ws.Range("G12").Value = "• JPMorgan Chase (JPMC) has elevated its Chief Data and Analytics Officer (CDAO) to a prominent Operating Committee role, signaling a firm-wide commitment to AI as a business imperative, not just a technology issue; 36% of companies now view AI as a business function reporting to senior leadership.`n• JPMC is deploying AI through an internal platform called LLM Suite, with over half of its employees now using it, demonstrating a cultural shift driven by both top-down strategy and bottom-up adoption; the firm focuses on AI applications that align with business objectives, such as fraud management, pricing, risk management, software engineering, and customer service in call centers.`n• JPMC recognizes data as a strategic asset and emphasizes the importance of a strong data foundation, including clean, understandable, and well-connected data across the organization, to enable successful AI implementation and drive business value; `"Having access to clean, understandable data is critical to our success.`""
$ws.Range("G13").Value = "•Commonwealth Bank of Australia (CBA) is cutting 45 jobs as part of a shift towards using AI for certain tasks.`n•The Finance Sector Union (FSU) claims a total of 90 roles are being eliminated, including 45 in direct banking, following the introduction of a new voice bot system for customer inquiries in June.`n•CBA is investing over A`$2 billion in operations, including technology, and states it is creating new roles and reskilling opportunities, while denying offshoring jobs."
